# Update comment strings: normalize delimiter examples from (', ')/('; ') to (", ")
# and fix grouping_columns / low_var_feature_removal default values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C5").Value = 'Patterns present in the names of pre-processing channels that will be used to remove them. Input separated by comma and whitespace (", "). This can be adjusted in the panel.csv file afterwards.'
$ws.Range("C8").Value = 'What ID your reference replicates have. Accepts multiple IDs, separated by comma and whitespace (", ").'
$ws.Range("B14").Value = 'antigen, anca_status'
$ws.Range("C14").Value = 'Columns in meta-table that contain group stratification. Accepts multiple names, separated by comma and whitespace (", ").'
$ws.Range("C15").Value = 'Control order of groups. Input group names separated by comma and whitespace (", "). For multiple grouping columns - separate by semicolon and whitespace (", "). If you do not want to input order for any specific grouping column - write NA.'
$ws.Range("C16").Value = 'Pre-gated subsets of events that you want to analyze. These should be the exact names of folders in fcs/4_subsets/. Separated by comma and whitespace (", ").'
$ws.Range("C17").Value = 'Whether to down- or upsample events PER SAMPLE using a given factor. You might want to downsample if studying e.g. granulocytes (on whole blood) to reduce computational burden and oversample if interested in finding rare populations. Downsampling could also be used to reduce the dominance of the oversized samples, while oversampling could be used to enhance the representation of the undersized ones. The cutoff for downsampling is the average sample size across all samples (samples are not reduced in size beyond that threshold). The cutoff for oversampling is the average sample size across all samples (samples are not increased in size beyond that threshold). Accepts multiple values separated by comma and whitespace (", ") if different needed for each data_subset. Order is given by data_subset order.'
$ws.Range("B19").Value = '0, 30'
$ws.Range("C19").Value = '0 (off) or 1 (on). Second value can be entered for number of features to keep, separated by comma and whitespace (", "). Is done by default for top 20 features if subset_feature_selection.xlsx is not provided.'
$ws.Range("C26").Value = 'If fs, number of clusters after hierarchical and before automatic merging. If fs_manual, final number of clusters after hierarchical merging. If pg, number of nearest neighbors parameter. Accepts multiple values, separated by comma and whitespace (", "), if different values are needed for each data_subset. Order is given by data_subset vector.'
$ws.Range("C30").Value = 'UMAPs n_neighbors. Low value leads to better resolution of local structure, high - to better capture of global structure. Default value is 15. Accepts multiple values, separated by comma and whitespace (", ") of different values are needed for each data_subset. Order is given by data_subset vector.'
$ws.Range("C31").Value = 'UMAPs min_dist. Lower values lead to clumpier embeddings, higher make them more spread-out. Accepts multiple values, separated by comma and whitespace (", ") of different values are needed for each data_subset. Order is given by data_subset vector.'

# Restore the scroll position / selection that was active when the author saved
$ws.Range("C17").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
